$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.187.08'
$ws.Range("E2").Value = '  -4.20%  '
$ws.Range("D3").Value = '1.658.35'
$ws.Range("E3").Value = '  -2.80%  '
$ws.Range("E4").Value = '  +0.22%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '217.15'
$ws.Range("E5").Value = '  -3.03%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5152'
$ws.Range("E6").Value = '  -3.25%  '
$ws.Range("E7").Value = '  +0.25%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2580'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06430'
$ws.Range("E9").Value = '  -2.68%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.98'
$ws.Range("E10").Value = '  -3.89%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07771'
$ws.Range("E11").Value = '  +1.55%  '
$ws.Range("D12").Value = '1.662.28'
$ws.Range("E12").Value = '  -2.67%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.291'
$ws.Range("E13").Value = '  -4.85%  '
$ws.Range("D14").Value = '1.887.19'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5536'
$ws.Range("E15").Value = '  -4.90%  '
$ws.Range("D16").Value = '0.0₅8043'
$ws.Range("E16").Value = '  -1.74%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.23'
$ws.Range("E17").Value = '  -5.09%  '
$ws.Range("D18").Value = '26.226.21'
$ws.Range("E18").Value = '  -4.08%  '
$ws.Range("E19").Value = '  +0.19%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '210.45'
$ws.Range("E20").Value = '  -2.57%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.418'
$ws.Range("E21").Value = '  -4.72%  '
$ws.Range("E22").Value = '  -3.29%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.032'
$ws.Range("E23").Value = '  +0.71%  '
$ws.Range("E24").Value = '  +0.24%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.34'
$ws.Range("E25").Value = '  +1.10%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.753'
$ws.Range("E26").Value = '  +2.98%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1168'
$ws.Range("E27").Value = '  -2.97%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.985'
$ws.Range("E28").Value = '  -3.35%  '
$ws.Range("E29").Value = '  -2.59%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05190'
$ws.Range("E30").Value = '  -3.43%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.251'
$ws.Range("E31").Value = '  -2.81%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.349'
$ws.Range("E32").Value = '  -3.82%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.228'
$ws.Range("E33").Value = '  -5.39%  '
$ws.Range("E34").Value = '  -4.42%  '
$ws.Range("E35").Value = '  -3.64%  '
$ws.Range("E36").Value = '  -1.27%  '
$ws.Range("E37").Value = '  -2.13%  '
$ws.Range("D38").Value = '1.178.22'
$ws.Range("E38").Value = '  +12.81%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5697'
$ws.Range("E39").Value = '  -2.74%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01593'
$ws.Range("E40").Value = '  -2.89%  '
$ws.Range("E41").Value = '  +0.18%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8391'
$ws.Range("E42").Value = '  -0.46%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.677'
$ws.Range("E43").Value = '  -2.36%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '100.34'
$ws.Range("E44").Value = '  -0.54%  '
$ws.Range("D45").Value = '1.797.42'
$ws.Range("E45").Value = '  -2.78%  '
$ws.Range("E46").Value = '  +4.87%  '
$ws.Range("E47").Value = '  +0.36%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '55.85'
$ws.Range("E48").Value = '  -3.51%  '
$ws.Range("E49").Value = '  -0.14%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.908'
$ws.Range("E50").Value = '  -2.18%  '
$ws.Range("E51").Value = '  -3.28%  '
